# Update the worksheet date and all the division-problem answers in the
# table to the new values for the 2026-02-11 edition.
#
# NOTE on ordering: the replacement list contains a value ("72÷3=24, 0")
# that is both an old value (row 2, col 4) and a new value (row 5, col 1).
# The replacements below are applied in the same order as they appear in
# the document (top-to-bottom, left-to-right within each row), so the
# original "72÷3=24, 0" cell is replaced (with "69÷9=7, 6") before the new
# "72÷3=24, 0" is written into a different cell later on. Because every
# other old value is unique in the document, exact whole-text matching
# via Find/Replace cannot cross-match the wrong cell.

$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-10 Tuesday", "2026-02-11 Wednesday"),
    @("92÷9=10, 2", "82÷9=9, 1"),
    @("48÷5=9, 3", "49÷4=12, 1"),
    @("95÷9=10, 5", "49÷6=8, 1"),
    @("23÷5=4, 3", "37÷4=9, 1"),
    @("75÷8=9, 3", "99÷2=49, 1"),
    @("19÷6=3, 1", "87÷3=29, 0"),
    @("73÷5=14, 3", "38÷7=5, 3"),
    @("14÷5=2, 4", "87÷9=9, 6"),
    @("72÷3=24, 0", "69÷9=7, 6"),
    @("63÷3=21, 0", "92÷6=15, 2"),
    @("71÷2=35, 1", "49÷2=24, 1"),
    @("78÷6=13, 0", "36÷9=4, 0"),
    @("22÷9=2, 4", "86÷3=28, 2"),
    @("45÷7=6, 3", "15÷2=7, 1"),
    @("16÷2=8, 0", "38÷6=6, 2"),
    @("48÷9=5, 3", "12÷2=6, 0"),
    @("41÷9=4, 5", "69÷2=34, 1"),
    @("38÷9=4, 2", "26÷5=5, 1"),
    @("95÷8=11, 7", "68÷6=11, 2"),
    @("53÷6=8, 5", "27÷3=9, 0"),
    @("26÷4=6, 2", "72÷3=24, 0"),
    @("17÷6=2, 5", "49÷3=16, 1"),
    @("19÷4=4, 3", "83÷2=41, 1"),
    @("30÷5=6, 0", "68÷9=7, 5"),
    @("22÷2=11, 0", "64÷4=16, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) {
        Write-Output "OK: '$old' -> '$new'"
    } else {
        Write-Output "MISSING: '$old' was not found (expected -> '$new')"
    }
}
